# AFDP-6132 fix user/groups names with correct prefix and domain
#
# The "Assignment Rules" sheet (Sheet1) contains several literal participant
# values of the form "owning group, <LDAP_ID>". Two distinct LDAP ids were
# used with the wrong prefix/domain:
#   ARKCASE_SUPERVISOR@ARMEDIA.COM              -> 000.ARKCASE_SUPERVISOR@APPDEV.ARMEDIA.COM
#   ARKCASE_ENTITY_ADMINISTRATOR@ARMEDIA.COM    -> 000.ARKCASE_ENTITY_ADMINISTRATOR@APPDEV.ARMEDIA.COM
#
# These values occur in column G of rows 23, 24, 30 (supervisor) and
# rows 34, 37 (entity administrator).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the "entity administrator" rows first, then the "supervisor" rows,
# so that new shared-string entries are introduced in that order.
$ws.Range("G34").Value = "owning group, 000.ARKCASE_ENTITY_ADMINISTRATOR@APPDEV.ARMEDIA.COM"
$ws.Range("G37").Value = "owning group, 000.ARKCASE_ENTITY_ADMINISTRATOR@APPDEV.ARMEDIA.COM"

$ws.Range("G23").Value = "owning group, 000.ARKCASE_SUPERVISOR@APPDEV.ARMEDIA.COM"
$ws.Range("G24").Value = "owning group, 000.ARKCASE_SUPERVISOR@APPDEV.ARMEDIA.COM"
$ws.Range("G30").Value = "owning group, 000.ARKCASE_SUPERVISOR@APPDEV.ARMEDIA.COM"

# Leave the view focused on the last-edited cell, matching the author's
# final selection when the workbook was saved.
$ws.Activate()
$ws.Range("G24").Select()
